# Add a new "canonical SMILES" column (D) to the Microstate List sheet,
# populating it with a de-stereo'd (slash/backslash-stripped) version of the
# "canonical isomeric SMILES" column (C) for every microstate row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Cells.Item(2, 4).Value = "canonical SMILES"

# microstate ID -> canonical SMILES (stereo bond markers removed from column C)
$canonicalSmiles = @{
    "SM03_micro001" = "c1ccc(cc1)Cc2[nH+]nc(s2)NC(=O)c3cccs3"
    "SM03_micro002" = "c1ccc(cc1)Cc2nnc(s2)N=C(c3cccs3)[O-]"
    "SM03_micro004" = "c1ccc(cc1)Cc2[nH+]nc(s2)N=C(c3cccs3)[O-]"
    "SM03_micro009" = "c1ccc(cc1)Cc2nnc(s2)[N-]C(=[OH+])c3cccs3"
    "SM03_micro010" = "c1ccc(cc1)[CH-]c2nnc(s2)NC(=O)c3cccs3"
    "SM03_micro011" = "c1ccc(cc1)Cc2nnc(s2)NC(=[OH+])c3cccs3"
    "SM03_micro012" = "c1ccc(cc1)Cc2nnc(s2)NC(=O)c3cccs3"
    "SM03_micro013" = "c1ccc(cc1)Cc2n[nH]c(=NC(=O)c3cccs3)s2"
    "SM03_micro014" = "c1ccc(cc1)[CH-]c2nnc(s2)N=C(c3cccs3)[O-]"
    "SM03_micro015" = "c1ccc(cc1)Cc2[nH+]nc(s2)[N-]C(=[OH+])c3cccs3"
    "SM03_micro016" = "c1ccc(cc1)Cc2n[nH+]c(s2)NC(=O)c3cccs3"
    "SM03_micro020" = "c1ccc(cc1)CC2=N[NH2+]C(=NC(=O)c3cccs3)S2"
    "SM03_micro021" = "c1ccc(cc1)CC2=[NH2+2]N=C(S2)NC(=O)c3cccs3"
    "SM03_micro022" = "c1ccc(cc1)Cc2n[nH]c(=NC(=[OH+])c3cccs3)s2"
    "SM03_micro023" = "c1ccc(cc1)Cc2[nH+][nH]c(=NC(=O)c3cccs3)s2"
    "SM03_micro024" = "c1ccc(cc1)Cc2[nH+][nH+]c(s2)NC(=O)c3cccs3"
}

for ($r = 3; $r -le 18; $r++) {
    $id = $ws.Cells.Item($r, 2).Value()
    $smiles = $canonicalSmiles[$id]
    if ($smiles) {
        $ws.Cells.Item($r, 4).Value = $smiles
    }
}

# Match the new column's width to the author's commit (closest value the
# engine's pixel-rounded column-width model can represent for a target of
# 36.85546875 "characters")
$ws.Range("D1:D18").ColumnWidth = 36.0
